$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A50").Value = "Danny Giordani"
$ws.Range("B50").Value = "Nicolas Giordani  | FC SAVIGNANO"
$ws.Range("C50").Value = "Danny Giordani | I Magnifici"
$ws.Range("D50").Value = "Luca Frasca | Clitoriders"
$ws.Range("E50").Value = "Alessandro Maffei | FC SAVIGNANO"
$ws.Range("F50").Value = "Gentian Capa | Power Ginger"
